$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data
$ws.Range('D2').Value = '70.419.78'
$ws.Range('E2').Value = '  +0.85%  '
$ws.Range('D3').Value = '3.622.09'
$ws.Range('E3').Value = '  +2.63%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '''601.97'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').Value = '''196.46'
$ws.Range('E6').Value = '  -0.12%  '
$ws.Range('E7').Value = '  -0.28%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  +6.14%  '
$ws.Range('E10').Value = '  -0.12%  '
$ws.Range('E11').Value = '  -0.83%  '
$ws.Range('E12').Value = '  +0.49%  '
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').Value = '4.198.46'
$ws.Range('D15').Value = '''606.87'
$ws.Range('E15').Value = '  +1.49%  '
$ws.Range('D16').Value = '''12.91'
$ws.Range('E16').Value = '  +1.73%  '
$ws.Range('D17').Value = '70.481.10'
$ws.Range('E17').Value = '  +0.74%  '
$ws.Range('D18').Value = '3.621.24'
$ws.Range('E18').Value = '  +2.99%  '
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').Value = '  +1.33%  '
$ws.Range('E21').Value = '  +0.94%  '
$ws.Range('D22').Value = '''18.22'
$ws.Range('E22').Value = '  -0.86%  '
$ws.Range('D23').Value = '''5.19'
$ws.Range('E23').Value = '  -1.65%  '
$ws.Range('D24').Value = '''103.10'
$ws.Range('E24').Value = '  +1.43%  '
$ws.Range('E25').Value = '  -0.06%  '
$ws.Range('E26').Value = '  -5.73%  '
$ws.Range('E27').Value = '  -2.22%  '
$ws.Range('E28').Value = '  +1.20%  '
$ws.Range('D29').Value = '''33.82'
$ws.Range('E29').Value = '  +1.35%  '
$ws.Range('D30').Value = '''4.70'
$ws.Range('E30').Value = '  +8.25%  '
$ws.Range('D31').Value = '''7.32'
$ws.Range('E31').Value = '  +3.57%  '
$ws.Range('D32').Value = '''12.28'
$ws.Range('E32').Value = '  -1.61%  '
$ws.Range('E33').Value = '  +2.00%  '
$ws.Range('D34').Value = '''63.49'
$ws.Range('E34').Value = '  +0.60%  '
$ws.Range('E35').Value = '  +2.98%  '
$ws.Range('D36').Value = '3.904.29'
$ws.Range('E36').Value = '  +4.75%  '
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').Value = '''3.07'
$ws.Range('E38').Value = '  -0.03%  '
$ws.Range('B39').Value = 'Bittensor'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D39').Value = '''517.27'
$ws.Range('E39').Value = '  +5.91%  '
$ws.Range('D40').Value = '''36.92'
$ws.Range('E41').Value = '  -0.84%  '
$ws.Range('E42').Value = '  -2.68%  '
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('D44').Value = '''0.0460'
$ws.Range('E44').Value = '  +1.40%  '
$ws.Range('D45').Value = '''3.55'
$ws.Range('E45').Value = '  +7.40%  '
$ws.Range('D46').Value = '''2.91'
$ws.Range('E46').Value = '  +2.93%  '
$ws.Range('E47').Value = '  -0.03%  '
$ws.Range('E48').Value = '  +0.35%  '
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('E50').Value = '  -0.55%  '
$ws.Range('E51').Value = '  +1.60%  '
